$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Row 11: test_ref_date_before_start ---
$ws1.Range("A11").Value = "test_ref_date_before_start"
$ws1.Range("C11").Value = "numpy.random"
$ws1.Range("D11").Value = "choice"
$ws1.Range("E11").Value = 1
$ws1.Range("F11").Value = 2
$ws1.Range("H11").Value = "kg"
$ws1.Range("H11").Font.Color = 0
$ws1.Range("I11").Value = 39814
$ws1.Range("I11").NumberFormat = "m/d/yy"
$ws1.Range("J11").Value = 39904
$ws1.Range("J11").NumberFormat = "m/d/yy"
$ws1.Range("K11").Value = 0.1
$ws1.Range("K11").NumberFormat = "0.00"
$ws1.Range("L11").Value = 39448
$ws1.Range("L11").NumberFormat = "m/d/yy"
$ws1.Range("M11").Value = "test var 1"

# --- Row 12: test_ref_date_after_end ---
$ws1.Range("A12").Value = "test_ref_date_after_end"
$ws1.Range("C12").Value = "numpy.random"
$ws1.Range("D12").Value = "choice"
$ws1.Range("E12").Value = 1
$ws1.Range("F12").Value = 2
$ws1.Range("H12").Value = "kg"
$ws1.Range("H12").Font.Color = 0
$ws1.Range("I12").Value = 39814
$ws1.Range("I12").NumberFormat = "m/d/yy"
$ws1.Range("J12").Value = 39904
$ws1.Range("J12").NumberFormat = "m/d/yy"
$ws1.Range("K12").Value = 0.1
$ws1.Range("K12").NumberFormat = "0.00"
$ws1.Range("L12").Value = 39448
$ws1.Range("L12").NumberFormat = "m/d/yy"
$ws1.Range("M12").Value = "test var 1"

# --- Row 13: variable with space at end in excel  ---
$ws1.Range("A13").Value = "variable with space at end in excel "
$ws1.Range("C13").Value = "numpy.random"
$ws1.Range("C13").Font.Color = 0
$ws1.Range("D13").Value = "choice"
$ws1.Range("D13").Font.Color = 0
$ws1.Range("E13").Value = 2
$ws1.Range("E13").Font.Color = 0
$ws1.Range("F13").Value = 4
$ws1.Range("F13").Font.Color = 0
$ws1.Range("G13").Font.Color = 0
$ws1.Range("H13").Value = "-"
$ws1.Range("H13").Font.Color = 0
$ws1.Range("M13").Value = "label"
$ws1.Range("M13").Font.Color = 0

# --- Row 14: empty row anchor ---
$ws1.Range("A14").Value = ""

$ws1.Range("C15").Select()
